# Actualización automática del inventario: agrega el nuevo producto
# "Film fusor inferior HP/CANON" (código EUG3NX) como fila 67.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 67

$ws.Cells.Item($row, 1).Value = "EUG3NX"
$ws.Cells.Item($row, 2).Value = "Film fusor inferior HP/CANON"
$ws.Cells.Item($row, 3).Value = "HP M252 M277 M377 M477 M452 M454 M455 M479, Canon MF 732 734 735"
$ws.Cells.Item($row, 4).Value = 40000
$ws.Cells.Item($row, 5).Value = 150000
$ws.Cells.Item($row, 6).Value = 3
$ws.Cells.Item($row, 7).Value = 4
$ws.Cells.Item($row, 8).Formula = "=(E67-D67)*G67"
$ws.Cells.Item($row, 9).Formula = "=D67*F67"
$ws.Cells.Item($row, 10).Value = 120000
